# Weekly fruta/hortaliza update: insert two new rows (Primera/Segunda) for a
# new reporting date (serial 44890) right after the current "first" week
# (row 397, serial 44225), pushing the rest of the historical rows down by
# two rows. Matches the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 396 (existing rows 396.. shift down to 398..)
$ws.Rows.Item(396).Insert()
$ws.Rows.Item(396).Insert()

# New row 396 - Calidad "Primera"
$ws.Range("A396").Value = 1
$ws.Range("B396").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C396").Value = "Arica y Parinacota"
$ws.Range("D396").Value = 44890
$ws.Range("E396").Value = 15
$ws.Range("F396").Value = 100114014
$ws.Range("G396").Value = "Betarraga"
$ws.Range("H396").Value = "Sin especificar"
$ws.Range("I396").Value = "Primera"
$ws.Range("J396").Value = 1000
$ws.Range("K396").Value = 400
$ws.Range("L396").Value = 450
$ws.Range("M396").Value = 430
$ws.Range("N396").Value = "$/paquete 4 unidades"
$ws.Range("O396").Value = "Región de Arica y Parinacota"
$ws.Range("P396").Value = 108
$ws.Range("Q396").Value = 4
$ws.Range("R396").Value = "Hortaliza"

# New row 397 - Calidad "Segunda"
$ws.Range("A397").Value = 1
$ws.Range("B397").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C397").Value = "Arica y Parinacota"
$ws.Range("D397").Value = 44890
$ws.Range("E397").Value = 15
$ws.Range("F397").Value = 100114014
$ws.Range("G397").Value = "Betarraga"
$ws.Range("H397").Value = "Sin especificar"
$ws.Range("I397").Value = "Segunda"
$ws.Range("J397").Value = 750
$ws.Range("K397").Value = 400
$ws.Range("L397").Value = 450
$ws.Range("M397").Value = 427
$ws.Range("N397").Value = "$/paquete 5 unidades"
$ws.Range("O397").Value = "Región de Arica y Parinacota"
$ws.Range("P397").Value = 85
$ws.Range("Q397").Value = 5
$ws.Range("R397").Value = "Hortaliza"

Write-Output "Inserted weekly rows 396-397; sheet now spans to row 407"
